$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sample-name labels in column A (rows 5-10) to new naming scheme.
# Original: t1, t2, t3, l1, l2, l3 (rows 5-10)
# New:      100nM1, 100nM3, 100nM2, 10uM1, 10uM2, 10uM3 (rows 5-10)
$ws.Range("A5").Value = "100nM1"
$ws.Range("A6").Value = "100nM3"
$ws.Range("A7").Value = "100nM2"
$ws.Range("A8").Value = "10uM1"
$ws.Range("A9").Value = "10uM2"
$ws.Range("A10").Value = "10uM3"
